$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G9").Value = 42
$ws.Range("G10").Value = 42
$ws.Range("G11").Value = 42
$ws.Range("G12").Value = 42
